# "Reference Done Chapter 1 Updated"
#
# The _GoBack bookmark (Word's "last edit location" marker) moves from the
# "Shortest Path Analysis" paragraph ("The game applies this algorithm...")
# to the end of the first sentence of the "Augmented Reality" paragraph
# ("...GPS data."), splitting that sentence's run into two runs the way
# Word does when a bookmark boundary falls inside an existing run.

$d = $word.ActiveDocument

# Remove the _GoBack bookmark from its old location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the first sentence of the Augmented Reality paragraph.
$target = $d.Content
$found = $target.Find.Execute(
    "Augmented reality is a live direct or indirect view of a physical, real-world environment whose elements are augmented (or supplemented) by computer-generated sensory input such as sound, video, graphics or GPS data.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Re-create _GoBack around that sentence (Word's implicit last-edit marker).
if ($found) {
    $d.Bookmarks.Add("_GoBack", $target)
}
